$d = $word.ActiveDocument

# --- 1. Small in-place wording tweaks (hyphenate frontend/backend) -------------
$d.Content.Find.Execute(
    "frontend backend developers",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "front-end back-end developers", 2) | Out-Null

# --- 2. "MageCloud Development Services" -> "MageCloud development Services" --
$d.Content.Find.Execute(
    "MageCloud Development Services",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "MageCloud development Services", 2) | Out-Null

# --- 3. Rewritten intro paragraph (after "MageCloud development Services") ----
$d.Paragraphs.Item(7).Range.Text = "Hello and welcome to MageCloud! We’re a top-rated, full-stack agency with an international client base. We’re different from other development providers in that we are committed to our clients’ continued success. We see the big picture of how your business could grow rather than remain focused on just one project. While other agencies might create a development site for you and leave all of the maintenance and upkeep to you and your team, at MageCloud, we offer continuous support for our customers, with the specific intention to continue to grow your Denmark eCommerce business. By using monthly analytics, we first propose and then implement solutions customized to fit your business’s needs so your profits continue to expand."

# --- 4. Paragraph under "Back-end Development" heading now holds the (updated)
#        extensions copy that used to sit under "Front-end Development" -------
$d.Paragraphs.Item(9).Range.Text = "Our developers can build new extensions from the ground up, or we can alter existing extensions depending on your business’s specific needs. Our extensions are tested via the latest stable version and all current web browsers before delivering or integrating to your development store so you never have to worry about a malfunction or slowing your website down."

# --- 5. Paragraph under "Front-end Development" heading now holds the new
#        front-end copy ---------------------------------------------------------
$d.Paragraphs.Item(11).Range.Text = "Front-end development is what your customers see when they’re on your website, so it should be both aesthetically interesting and user-friendly, no matter what device your customers are logging on from. Our developers build mobile-friendly sites that are specifically designed for cell phones, integrate third party widgets and extensions, and offer the most up-to-date front-end customization."

# --- 6. "Platform Integrations" block moves up to directly follow the
#        front-end copy; its body paragraph is the (updated) old back-end copy -
$d.Paragraphs.Item(12).Range.Text = "Platform <strong>Integrations</strong>"
$d.Paragraphs.Item(13).Range.Text = "Back-end development might sound confusing, but in the most simple of terms, it is the portion of your eCommerce site that you log into in order to make changes or update inventory. We use back-end development to help our clients manage their development stores. Our team ensures that our clients can maintain their physical inventory across multiple warehouses and showrooms, tracking any movement of items from one location to another. We’re able to give our clients the ability to segment the price creation process into three basic elements: defining pricing rules, assigning rules to create specific pricing actions, and designating pricing controls. Finally, we assist clients in creating multiple stores, each store is then accessible on different URLs, under the same installation using a shared shopping cart. These share the same backend, including sharing an inventory,which makes the administration aspect simple and easy for our clients."

# --- 7. "Plugin Development" block follows, with reworded integrations copy ---
$d.Paragraphs.Item(14).Range.Text = "Plugin <strong>Development</strong>"
$d.Paragraphs.Item(15).Range.Text = "The more integrations your web platform supports, the smoother your development store will run and the less time you will need to spend working on your website, freeing up your time to dedicate to other areas of your business that need your attention. Some of the platform integrations we offer are: shipping platforms integrations, payment platforms integration, inventory management integrations, and Salesforce integration."

# --- 8. "Why Pick MageCloud?" now comes last, with brand-new closing copy -----
$d.Paragraphs.Item(16).Range.Text = "Why Pick MageCloud?"
$d.Paragraphs.Item(17).Range.Text = "We’re not like other developers. We always maintain consistent and reliable communication with our clients through the entire development process to ensure our goals are completely aligned and they can anticipate what the face of their business will look like when their development site goes live. Depending on our client’s goals and objectives, we develop realistic plans that meet their every expectation so they can continue to grow their business."
